$wb = $excel.ActiveWorkbook

# --- Sheet "Y" (sheet1.xml): add new shock row data, clear old A3:A6 values ---
$wsY = $wb.Worksheets.Item("Y")
$wsY.Range("B2").Value = "Maize (home consumed)"
$wsY.Range("B2").Style = "Normal"
$wsY.Range("C2").Value = 100000000000000
$wsY.Range("A3").Value = $null
$wsY.Range("A4").Value = $null
$wsY.Range("A5").Value = $null
$wsY.Range("A6").Value = $null

# --- Sheet "A" (sheet2.xml): move selection, drop frozen/scrolled view state ---
$wsA = $wb.Worksheets.Item("A")
$wsA.Range("F2").Select()

# --- Sheet "VA" (sheet3.xml): move selection ---
$wsVA = $wb.Worksheets.Item("VA")
$wsVA.Range("D10").Select()

# --- Sheet "Y" becomes the active tab/selection last, so it ends up the selected tab ---
$wsY.Activate()
$wsY.Range("B3").Select()
